# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the Famfrit Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 10825.3
$ws.Range("I9").Value = 16766.834
$ws.Range("J9").Value = 1913
$ws.Range("K9").Value = 16766.834
$ws.Range("L9").Value = 1913
$ws.Range("M9").Value = -16597.834
$ws.Range("N9").Value = -2251
$ws.Range("H17").Value = 4091558.2
$ws.Range("J17").Value = 4091558.2
$ws.Range("L17").Value = 12274674.6
$ws.Range("N17").Value = -12275010.6
$ws.Range("H38").Value = 3228.1562
$ws.Range("J38").Value = 3914
$ws.Range("L38").Value = 11742
$ws.Range("N38").Value = -12486
$ws.Range("H106").Value = 2749
$ws.Range("I106").Value = 2749
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2749
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2118
$ws.Range("H113").Value = 5160
$ws.Range("I113").Value = 5380.6665
$ws.Range("K113").Value = 5380.6665
$ws.Range("M113").Value = -2126.6665
$ws.Range("H121").Value = 10000
$ws.Range("J121").Value = 10000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494
$ws.Range("H132").Value = 3375.3865
$ws.Range("I132").Value = 3418.5642
$ws.Range("K132").Value = 10255.6926
$ws.Range("M132").Value = -7725.692599999998
$ws.Range("H137").Value = 11334.775
$ws.Range("J137").Value = 9558.799999999999
$ws.Range("L137").Value = 28676.4
$ws.Range("N137").Value = -33776.39999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4000.8572
$ws.Range("I2").Value = 2749.182
$ws.Range("K2").Value = 2749.182
$ws.Range("M2").Value = -2636.182
$ws.Range("H32").Value = 4822.073
$ws.Range("I32").Value = 4505.125
$ws.Range("K32").Value = 4505.125
$ws.Range("M32").Value = -4218.125
$ws.Range("H45").Value = 2602.4119
$ws.Range("I45").Value = 2061.25
$ws.Range("J45").Value = 3375.5
$ws.Range("K45").Value = 2061.25
$ws.Range("L45").Value = 3375.5
$ws.Range("M45").Value = -1684.25
$ws.Range("N45").Value = -4129.5
$ws.Range("H61").Value = 3935.625
$ws.Range("I61").Value = 3000
$ws.Range("J61").Value = 4247.5
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 4247.5
$ws.Range("M61").Value = -2788
$ws.Range("N61").Value = -4671.5
$ws.Range("H116").Value = 4000.8572
$ws.Range("I116").Value = 2749.182
$ws.Range("K116").Value = 2749.182
$ws.Range("M116").Value = -455.1819999999998
$ws.Range("H132").Value = 47245.03
$ws.Range("I132").Value = 5257.231
$ws.Range("K132").Value = 15771.693
$ws.Range("M132").Value = -13241.693
$ws.Range("H136").Value = 3935.625
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 4247.5
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 12742.5
$ws.Range("M136").Value = -6450
$ws.Range("N136").Value = -17842.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4000.8572
$ws.Range("I3").Value = 2749.182
$ws.Range("K3").Value = 2749.182
$ws.Range("M3").Value = -2635.182
$ws.Range("H86").Value = 78062.5
$ws.Range("I86").Value = 63625
$ws.Range("K86").Value = 63625
$ws.Range("M86").Value = -62502
$ws.Range("H89").Value = 78062.5
$ws.Range("I89").Value = 63625
$ws.Range("K89").Value = 318125
$ws.Range("M89").Value = -312509
$ws.Range("H105").Value = 5287.625
$ws.Range("I105").Value = 3537.5386
$ws.Range("J105").Value = 6130.2593
$ws.Range("K105").Value = 3537.5386
$ws.Range("L105").Value = 6130.2593
$ws.Range("M105").Value = -1790.5386
$ws.Range("N105").Value = -9624.2593

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8336526
$ws.Range("I31").Value = 2385
$ws.Range("J31").Value = 31255414
$ws.Range("K31").Value = 2385
$ws.Range("L31").Value = 31255414
$ws.Range("M31").Value = -2090
$ws.Range("N31").Value = -31256004
$ws.Range("H34").Value = 8336526
$ws.Range("I34").Value = 2385
$ws.Range("J34").Value = 31255414
$ws.Range("K34").Value = 2385
$ws.Range("L34").Value = 31255414
$ws.Range("M34").Value = -2183
$ws.Range("N34").Value = -31255818
$ws.Range("H57").Value = 39499.5
$ws.Range("J57").Value = 39499.5
$ws.Range("L57").Value = 39499.5
$ws.Range("N57").Value = -40619.5
$ws.Range("H132").Value = 3984.2856
$ws.Range("I132").Value = 3840.5
$ws.Range("K132").Value = 11521.5
$ws.Range("M132").Value = -8991.5
$ws.Range("H134").Value = 3774.4092
$ws.Range("I134").Value = 3061.5881
$ws.Range("K134").Value = 9184.764299999999
$ws.Range("M134").Value = -6649.764299999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5776
$ws.Range("I3").Value = 1501.5714
$ws.Range("J3").Value = 15749.667
$ws.Range("K3").Value = 4504.7142
$ws.Range("L3").Value = 47249.001
$ws.Range("M3").Value = -4392.7142
$ws.Range("N3").Value = -47473.001
$ws.Range("H18").Value = 3891.25
$ws.Range("I18").Value = 1855
$ws.Range("K18").Value = 5565
$ws.Range("M18").Value = -5396
$ws.Range("H109").Value = 4041.4
$ws.Range("I109").Value = 8125.4
$ws.Range("J109").Value = 1999.4
$ws.Range("K109").Value = 24376.2
$ws.Range("L109").Value = 5998.200000000001
$ws.Range("M109").Value = -23336.2
$ws.Range("N109").Value = -8078.200000000001
$ws.Range("H131").Value = 35715824
$ws.Range("J131").Value = 1904.875
$ws.Range("L131").Value = 5714.625
$ws.Range("N131").Value = -15794.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1716.6923
$ws.Range("I97").Value = 756.4286
$ws.Range("K97").Value = 756.4286
$ws.Range("M97").Value = -260.4286
$ws.Range("H113").Value = 2933.3333
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2933.3333
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2933.3333
$ws.Range("N113").Value = -7273.3333
$ws.Range("H132").Value = 3373.6667
$ws.Range("I132").Value = 2560.6667
$ws.Range("K132").Value = 7682.000100000001
$ws.Range("M132").Value = -5152.000100000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5819.4443
$ws.Range("I7").Value = 5270.7144
$ws.Range("J7").Value = 7740
$ws.Range("K7").Value = 5270.7144
$ws.Range("L7").Value = 7740
$ws.Range("M7").Value = -5158.7144
$ws.Range("N7").Value = -7964
$ws.Range("H22").Value = 2154.077
$ws.Range("I22").Value = 1625.125
$ws.Range("J22").Value = 3000.4
$ws.Range("K22").Value = 1625.125
$ws.Range("L22").Value = 3000.4
$ws.Range("M22").Value = -1330.125
$ws.Range("N22").Value = -3590.4
$ws.Range("H27").Value = 2154.077
$ws.Range("I27").Value = 1625.125
$ws.Range("J27").Value = 3000.4
$ws.Range("K27").Value = 1625.125
$ws.Range("L27").Value = 3000.4
$ws.Range("M27").Value = -1518.125
$ws.Range("N27").Value = -3214.4
$ws.Range("H68").Value = 8459.4
$ws.Range("J68").Value = 7432.6665
$ws.Range("L68").Value = 7432.6665
$ws.Range("N68").Value = -8930.666499999999
$ws.Range("H71").Value = 8459.4
$ws.Range("J71").Value = 7432.6665
$ws.Range("L71").Value = 37163.3325
$ws.Range("N71").Value = -44651.3325
$ws.Range("H126").Value = 5819.4443
$ws.Range("I126").Value = 5270.7144
$ws.Range("J126").Value = 7740
$ws.Range("K126").Value = 15812.1432
$ws.Range("L126").Value = 23220
$ws.Range("M126").Value = -13342.1432
$ws.Range("N126").Value = -28160

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 28498
$ws.Range("J52").Value = 40495
$ws.Range("L52").Value = 40495
$ws.Range("N52").Value = -40947
$ws.Range("H70").Value = 27498.334
$ws.Range("J70").Value = 35000
$ws.Range("L70").Value = 35000
$ws.Range("N70").Value = -35630
$ws.Range("H73").Value = 27498.334
$ws.Range("J73").Value = 35000
$ws.Range("L73").Value = 35000
$ws.Range("N73").Value = -37184
$ws.Range("H96").Value = 7567.1665
$ws.Range("J96").Value = 7000.75
$ws.Range("L96").Value = 7000.75
$ws.Range("N96").Value = -9746.75
$ws.Range("H107").Value = 1082.6786
$ws.Range("I107").Value = 815.56525
$ws.Range("K107").Value = 2446.69575
$ws.Range("M107").Value = -526.6957499999999
$ws.Range("I132").Value = 2125.2173
$ws.Range("J132").Value = 4299.8
$ws.Range("K132").Value = 6375.651899999999
$ws.Range("L132").Value = 12899.4
$ws.Range("M132").Value = -3845.651899999999
$ws.Range("N132").Value = -17959.4
$ws.Range("H136").Value = 3593.0715
$ws.Range("I136").Value = 2333.5
$ws.Range("J136").Value = 11150.5
$ws.Range("K136").Value = 7000.5
$ws.Range("L136").Value = 33451.5
$ws.Range("M136").Value = -4450.5
$ws.Range("N136").Value = -38551.5

# ---- Cell deletions (values removed entirely by the refresh) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N106").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M113").ClearContents()
